$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = "35.437.91" }
    @{ Cell = "E2"; Value = "  +1.02%  " }
    @{ Cell = "D3"; Value = "1.903.49" }
    @{ Cell = "E3"; Value = "  +2.40%  " }
    @{ Cell = "E4"; Value = "  +0.26%  " }
    @{ Cell = "D5"; Value = "246.22" }
    @{ Cell = "E5"; Value = "  +3.64%  " }
    @{ Cell = "D6"; Value = "0.640" }
    @{ Cell = "E6"; Value = "  +2.95%  " }
    @{ Cell = "E7"; Value = "  +0.24%  " }
    @{ Cell = "D8"; Value = "41.90" }
    @{ Cell = "E8"; Value = "  -1.56%  " }
    @{ Cell = "D9"; Value = "0.340" }
    @{ Cell = "E9"; Value = "  +3.60%  " }
    @{ Cell = "E10"; Value = "  +1.33%  " }
    @{ Cell = "D11"; Value = "0.0999" }
    @{ Cell = "E11"; Value = "  +0.89%  " }
    @{ Cell = "D12"; Value = "2.180.18" }
    @{ Cell = "E12"; Value = "  +2.44%  " }
    @{ Cell = "D13"; Value = "12.39" }
    @{ Cell = "E13"; Value = "  +8.61%  " }
    @{ Cell = "E14"; Value = "  +2.81%  " }
    @{ Cell = "D15"; Value = "1.913.91" }
    @{ Cell = "E15"; Value = "  +2.97%  " }
    @{ Cell = "E16"; Value = "  +2.51%  " }
    @{ Cell = "D17"; Value = "35.484.56" }
    @{ Cell = "E17"; Value = "  +1.23%  " }
    @{ Cell = "D18"; Value = "71.96" }
    @{ Cell = "E18"; Value = "  +2.36%  " }
    @{ Cell = "D19"; Value = "0.0₃0829" }
    @{ Cell = "E19"; Value = "  +3.97%  " }
    @{ Cell = "D20"; Value = "243.09" }
    @{ Cell = "E20"; Value = "  +0.88%  " }
    @{ Cell = "D21"; Value = "12.66" }
    @{ Cell = "E21"; Value = "  +3.94%  " }
    @{ Cell = "E22"; Value = "  +1.56%  " }
    @{ Cell = "E23"; Value = "  +0.24%  " }
    @{ Cell = "E24"; Value = "  +1.04%  " }
    @{ Cell = "E25"; Value = "  +15.33%  " }
    @{ Cell = "D26"; Value = "171.70" }
    @{ Cell = "E26"; Value = "  +0.27%  " }
    @{ Cell = "E27"; Value = "  +7.75%  " }
    @{ Cell = "D28"; Value = "17.97" }
    @{ Cell = "E28"; Value = "  +1.67%  " }
    @{ Cell = "E29"; Value = "  +1.01%  " }
    @{ Cell = "D30"; Value = "0.978" }
    @{ Cell = "E30"; Value = "  +25.08%  " }
    @{ Cell = "E31"; Value = "  +1.96%  " }
    @{ Cell = "E32"; Value = "  +2.90%  " }
    @{ Cell = "E33"; Value = "  +0.22%  " }
    @{ Cell = "D34"; Value = "4.17" }
    @{ Cell = "E34"; Value = "  +4.03%  " }
    @{ Cell = "E35"; Value = "  +8.57%  " }
    @{ Cell = "E36"; Value = "  +0.02%  " }
    @{ Cell = "D37"; Value = "1.33" }
    @{ Cell = "E37"; Value = "  +2.76%  " }
    @{ Cell = "E38"; Value = "  +2.05%  " }
    @{ Cell = "D39"; Value = "0.0641" }
    @{ Cell = "E39"; Value = "  +17.16%  " }
    @{ Cell = "B40"; Value = "VeChain" }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "D40"; Value = "0.0204" }
    @{ Cell = "E40"; Value = "  +1.02%  " }
    @{ Cell = "B41"; Value = "Aave" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave" }
    @{ Cell = "D41"; Value = "91.68" }
    @{ Cell = "E41"; Value = "  -0.26%  " }
    @{ Cell = "B42"; Value = "MultiversX" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld" }
    @{ Cell = "D42"; Value = "50.53" }
    @{ Cell = "E42"; Value = "  +44.91%  " }
    @{ Cell = "B43"; Value = "InjectiveProtocol" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj" }
    @{ Cell = "D43"; Value = "15.59" }
    @{ Cell = "E43"; Value = "  +4.61%  " }
    @{ Cell = "B44"; Value = "Maker" }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" }
    @{ Cell = "D44"; Value = "1.347.02" }
    @{ Cell = "E44"; Value = "  -0.42%  " }
    @{ Cell = "E45"; Value = "  +2.15%  " }
    @{ Cell = "D46"; Value = "12.89" }
    @{ Cell = "E46"; Value = "  +2.21%  " }
    @{ Cell = "D47"; Value = "2.41" }
    @{ Cell = "E47"; Value = "  +0.08%  " }
    @{ Cell = "D48"; Value = "2.76" }
    @{ Cell = "E48"; Value = "  +0.00%  " }
    @{ Cell = "E49"; Value = "  +4.20%  " }
    @{ Cell = "D50"; Value = "2.091.43" }
    @{ Cell = "E50"; Value = "  +2.26%  " }
    @{ Cell = "D51"; Value = "0.0692" }
    @{ Cell = "E51"; Value = "  +1.66%  " }
)

foreach ($ch in $changes) {
    $r = $ws.Range($ch.Cell)
    $r.NumberFormat = "@"
    $r.Value = $ch.Value
    $r.Style = "Normal"
}

Write-Output "Applied $($changes.Count) cell updates"